$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = "28.302.52"
$c.Style = "Normal"
$ws.Range("E2").Value = "  +4.35%  "
$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = "1.803.33"
$c.Style = "Normal"
$ws.Range("E3").Value = "  +2.25%  "
$c = $ws.Range("D4")
$c.NumberFormat = "@"
$c.Value = "0.9993"
$c.Style = "Normal"
$ws.Range("E4").Value = "  +0.55%  "
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "316.29"
$c.Style = "Normal"
$ws.Range("E5").Value = "  +1.26%  "
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "0.9995"
$c.Style = "Normal"
$ws.Range("E6").Value = "  +0.37%  "
$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = "0.5532"
$c.Style = "Normal"
$ws.Range("E7").Value = "  +6.45%  "
$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = "0.3869"
$c.Style = "Normal"
$ws.Range("E8").Value = "  +7.29%  "
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = "0.07584"
$c.Style = "Normal"
$ws.Range("E9").Value = "  +4.08%  "
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "42.74"
$c.Style = "Normal"
$ws.Range("E10").Value = "  +1.04%  "
$ws.Range("E11").Value = "  +4.37%  "
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = "0.9995"
$c.Style = "Normal"
$ws.Range("E12").Value = "  +0.66%  "
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = "21.10"
$c.Style = "Normal"
$ws.Range("E13").Value = "  +3.17%  "
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "6.201"
$c.Style = "Normal"
$ws.Range("E14").Value = "  +3.01%  "
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "7.367"
$c.Style = "Normal"
$ws.Range("E15").Value = "  +7.14%  "
$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = "1.801.90"
$c.Style = "Normal"
$ws.Range("E16").Value = "  +3.31%  "
$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = "92.29"
$c.Style = "Normal"
$ws.Range("E17").Value = "  +5.56%  "
$ws.Range("E18").Value = "  +2.81%  "
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "0.06445"
$c.Style = "Normal"
$ws.Range("E19").Value = "  +0.42%  "
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "0.9994"
$c.Style = "Normal"
$ws.Range("E20").Value = "  +0.26%  "
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "17.32"
$c.Style = "Normal"
$ws.Range("E21").Value = "  +3.80%  "
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "5.987"
$c.Style = "Normal"
$ws.Range("E22").Value = "  +3.45%  "
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "28.319.66"
$c.Style = "Normal"
$ws.Range("E23").Value = "  +4.29%  "
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "11.44"
$c.Style = "Normal"
$ws.Range("E24").Value = "  +1.04%  "
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "2.133"
$c.Style = "Normal"
$ws.Range("E25").Value = "  +4.06%  "
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = "157.85"
$c.Style = "Normal"
$ws.Range("E26").Value = "  +3.11%  "
$ws.Range("E27").Value = "  +2.85%  "
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "2.402"
$c.Style = "Normal"
$ws.Range("E28").Value = "  +6.46%  "
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = "2.007.41"
$c.Style = "Normal"
$ws.Range("E29").Value = "  +3.28%  "
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = "123.84"
$c.Style = "Normal"
$ws.Range("E30").Value = "  +2.85%  "
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "1.122"
$c.Style = "Normal"
$ws.Range("E31").Value = "  +6.51%  "
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "0.1021"
$c.Style = "Normal"
$ws.Range("E32").Value = "  +6.54%  "
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "5.741"
$c.Style = "Normal"
$ws.Range("E33").Value = "  +5.61%  "
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "3.664"
$c.Style = "Normal"
$ws.Range("E34").Value = "  +1.97%  "
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = "0.2344"
$c.Style = "Normal"
$ws.Range("E35").Value = "  +16.79%  "
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "0.06289"
$c.Style = "Normal"
$ws.Range("E36").Value = "  +6.12%  "
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "8.901"
$c.Style = "Normal"
$ws.Range("E37").Value = "  +15.72%  "
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "0.02319"
$c.Style = "Normal"
$ws.Range("E38").Value = "  +5.32%  "
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = "11.63"
$c.Style = "Normal"
$ws.Range("E39").Value = "  +4.83%  "
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "5.047"
$c.Style = "Normal"
$ws.Range("E40").Value = "  +5.07%  "
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "0.6402"
$c.Style = "Normal"
$ws.Range("E41").Value = "  +5.16%  "
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "0.9991"
$c.Style = "Normal"
$ws.Range("E42").Value = "  +0.46%  "
$ws.Range("E43").Value = "  +3.80%  "
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "1.381"
$c.Style = "Normal"
$ws.Range("E44").Value = "  -3.09%  "
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "13.48"
$c.Style = "Normal"
$ws.Range("E45").Value = "  +3.77%  "
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "0.5977"
$c.Style = "Normal"
$ws.Range("E46").Value = "  +4.79%  "
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "3.684"
$c.Style = "Normal"
$ws.Range("E47").Value = "  +2.30%  "
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "124.65"
$c.Style = "Normal"
$ws.Range("E48").Value = "  +3.24%  "
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "1.974"
$c.Style = "Normal"
$ws.Range("E49").Value = "  +5.94%  "
$ws.Range("E50").Value = "  +3.95%  "
$ws.Range("E51").Value = "  +3.40%  "
